# Add two new Mac-Addresses (10 new device rows) to the registration center
# device master table, continuing the existing data pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$startMac = 3000166
$rowCount = 10

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $mac = $startMac + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $mac
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Reposition the view/selection to reflect where the user ended up after
# entering the new rows.
$excel.ActiveWindow.ScrollRow = 140
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E155").Select()
